# Add a "margin of victory" flag column, expressed as two new summary
# sheets: "margin_5_or_less" (margin_flag x stance) and
# "margin_5_or_less_withprez" (p16winningparty x margin_flag x stance).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert "margin_5_or_less" right after "rural_morethanfifth" (so it
#    becomes sheet #6, pushing the four "*_andprezresults" sheets down).
# ---------------------------------------------------------------------
$afterRural = $wb.Worksheets.Item("rural_morethanfifth")
$marginSheet = $wb.Worksheets.Add($null, $afterRural)
$marginSheet.Name = "margin_5_or_less"

$marginSheet.Range("A1").Value = "margin_flag"
$marginSheet.Range("B1").Value = "stance"
$marginSheet.Range("C1").Value = "n"
$headerRange = $marginSheet.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

$marginSheet.Range("A2").Value = "5_points_or_less"
$marginSheet.Range("B2").Value = "not_sponsoring"
$marginSheet.Range("C2").Value = 11

$marginSheet.Range("A3").Value = "5_points_or_less"
$marginSheet.Range("B3").Value = "sponsoring"
$marginSheet.Range("C3").Value = 15

$marginSheet.Range("A4").Value = "more_than_5_points"
$marginSheet.Range("B4").Value = "not_sponsoring"
$marginSheet.Range("C4").Value = 14

$marginSheet.Range("A5").Value = "more_than_5_points"
$marginSheet.Range("B5").Value = "sponsoring"
$marginSheet.Range("C5").Value = 190

$marginSheet.Range("A6").Value = "other"
$marginSheet.Range("B6").Value = "sponsoring"
$marginSheet.Range("C6").Value = 4

# ---------------------------------------------------------------------
# 2) Append "margin_5_or_less_withprez" as the new last sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$marginPrezSheet = $wb.Worksheets.Add($null, $lastSheet)
$marginPrezSheet.Name = "margin_5_or_less_withprez"

$marginPrezSheet.Range("A1").Value = "p16winningparty"
$marginPrezSheet.Range("B1").Value = "margin_flag"
$marginPrezSheet.Range("C1").Value = "stance"
$marginPrezSheet.Range("D1").Value = "n"
$headerRange2 = $marginPrezSheet.Range("A1:D1")
$headerRange2.Font.Bold = $true
$headerRange2.HorizontalAlignment = -4108

$marginPrezSheet.Range("A2").Value = "D"
$marginPrezSheet.Range("B2").Value = "5_points_or_less"
$marginPrezSheet.Range("C2").Value = "not_sponsoring"
$marginPrezSheet.Range("D2").Value = 1

$marginPrezSheet.Range("A3").Value = "D"
$marginPrezSheet.Range("B3").Value = "5_points_or_less"
$marginPrezSheet.Range("C3").Value = "sponsoring"
$marginPrezSheet.Range("D3").Value = 11

$marginPrezSheet.Range("A4").Value = "D"
$marginPrezSheet.Range("B4").Value = "more_than_5_points"
$marginPrezSheet.Range("C4").Value = "not_sponsoring"
$marginPrezSheet.Range("D4").Value = 9

$marginPrezSheet.Range("A5").Value = "D"
$marginPrezSheet.Range("B5").Value = "more_than_5_points"
$marginPrezSheet.Range("C5").Value = "sponsoring"
$marginPrezSheet.Range("D5").Value = 178

$marginPrezSheet.Range("A6").Value = "D"
$marginPrezSheet.Range("B6").Value = "other"
$marginPrezSheet.Range("C6").Value = "sponsoring"
$marginPrezSheet.Range("D6").Value = 4

$marginPrezSheet.Range("A7").Value = "R"
$marginPrezSheet.Range("B7").Value = "5_points_or_less"
$marginPrezSheet.Range("C7").Value = "not_sponsoring"
$marginPrezSheet.Range("D7").Value = 10

$marginPrezSheet.Range("A8").Value = "R"
$marginPrezSheet.Range("B8").Value = "5_points_or_less"
$marginPrezSheet.Range("C8").Value = "sponsoring"
$marginPrezSheet.Range("D8").Value = 4

$marginPrezSheet.Range("A9").Value = "R"
$marginPrezSheet.Range("B9").Value = "more_than_5_points"
$marginPrezSheet.Range("C9").Value = "not_sponsoring"
$marginPrezSheet.Range("D9").Value = 5

$marginPrezSheet.Range("A10").Value = "R"
$marginPrezSheet.Range("B10").Value = "more_than_5_points"
$marginPrezSheet.Range("C10").Value = "sponsoring"
$marginPrezSheet.Range("D10").Value = 12

# ---------------------------------------------------------------------
# 3) Restore per-sheet selection state as it was left after editing.
#    Sheets that were merely paged through settle on C22; the
#    originally active sheet (prezresults2016) ends up re-selected at
#    K7 and remains the active tab.
# ---------------------------------------------------------------------
[void]$wb.Worksheets.Item("gdp_vs_nationalavg").Range("C22").Select()
[void]$wb.Worksheets.Item("college_vs_nationalavg").Range("C22").Select()
[void]$wb.Worksheets.Item("nonwhite_vs_nationalavg").Range("C22").Select()
[void]$wb.Worksheets.Item("rural_morethanfifth").Range("C22").Select()
[void]$wb.Worksheets.Item("margin_5_or_less").Range("C22").Select()
[void]$wb.Worksheets.Item("gdp_andprezresults").Range("C22").Select()
[void]$wb.Worksheets.Item("college_degree_andprezresults").Range("C22").Select()
[void]$wb.Worksheets.Item("nonwhite_pop_andprezresults").Range("C22").Select()
[void]$wb.Worksheets.Item("rural_area_andprezresults").Range("C22").Select()

$prez = $wb.Worksheets.Item("prezresults2016")
[void]$prez.Activate()
[void]$prez.Range("K7").Select()
